$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = "a"
$ws.Range("D4").Value = "b"
$ws.Range("E4").Value = "Sup1"
$ws.Range("F4").Value = 12
$ws.Range("I4").Value = 32
$ws.Range("J4").Value = 21
$ws.Range("M4").Value = "Group1"

# Row 5
$ws.Range("C5").Value = "aa"
$ws.Range("D5").Value = "bb"
$ws.Range("E5").Value = "Sup2"
$ws.Range("F5").Value = 132
$ws.Range("I5").Value = 132
$ws.Range("J5").Value = 33
$ws.Range("M5").Value = "Group1"

# Row 6
$ws.Range("C6").Value = "aaa"
$ws.Range("D6").Value = "bbb"
$ws.Range("E6").Value = "Sup1"
$ws.Range("F6").Value = 32
$ws.Range("I6").Value = 42
$ws.Range("J6").Value = 12
$ws.Range("M6").Value = "Group2"

# Row 7
$ws.Range("C7").Value = "aaaa"
$ws.Range("D7").Value = "bbbb"
$ws.Range("E7").Value = "Sup1"
$ws.Range("F7").Value = 312
$ws.Range("I7").Value = 231
$ws.Range("J7").Value = 321
$ws.Range("M7").Value = "Group2"

# Rows 8-11: clear content (product removed from list)
$ws.Range("C8:F11").ClearContents()
$ws.Range("I8:J11").ClearContents()
$ws.Range("M8:M11").ClearContents()

# G column: enter 0 for rows 4-99
$gRange = $ws.Range("G4:G99")
$gRange.Value = 0
$gRange.NumberFormat = "0"

Write-Host "done"
